# Update countries & provincias Spain
# Applies the data refresh captured in the commit diff:
#  - Update "last updated" timestamp (A1)
#  - Update Kazajistan row values
#  - Belgica overtakes Kuwait in total cases -> rows swap, values updated
#  - Bahamas overtakes Benin in total cases -> rows swap, values updated

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 31 de Agosto de 2020 a las 05:53"

# 2. Row 33 - Kazajistan: refresh counters (label/country unchanged)
$ws.Cells.Item(33, 2).Value = 105795   # Casos totales
$ws.Cells.Item(33, 3).Value = 111      # Nuevos casos
$ws.Cells.Item(33, 4).Value = 96297    # Recuperados
$ws.Cells.Item(33, 5).Value = 7975     # Casos activos

# 3. Rows 40/41 - Belgica now has more total cases than Kuwait, so it
#    takes row 40 (with refreshed data) and Kuwait drops to row 41
#    (values unchanged from before).
$ws.Cells.Item(40, 1).Value = "Belgica"
$ws.Cells.Item(40, 2).Value = 85042
$ws.Cells.Item(40, 3).Value = 443
$ws.Cells.Item(40, 4).Value = 18415
$ws.Cells.Item(40, 5).Value = 56733
$ws.Cells.Item(40, 6).Value = 0
$ws.Cells.Item(40, 7).Value = 3
$ws.Cells.Item(40, 8).Value = 9894

$ws.Cells.Item(41, 1).Value = "Kuwait"
$ws.Cells.Item(41, 2).Value = 84636
$ws.Cells.Item(41, 3).Value = 0
$ws.Cells.Item(41, 4).Value = 76650
$ws.Cells.Item(41, 5).Value = 7456
$ws.Cells.Item(41, 6).Value = 0
$ws.Cells.Item(41, 7).Value = 0
$ws.Cells.Item(41, 8).Value = 530

# 4. Rows 138/139 - Bahamas now has more total cases than Benin, so it
#    takes row 138 (with refreshed data) and Benin drops to row 139
#    (values unchanged from before).
$ws.Cells.Item(138, 1).Value = "Bahamas"
$ws.Cells.Item(138, 2).Value = 2167
$ws.Cells.Item(138, 3).Value = 0
$ws.Cells.Item(138, 4).Value = 782
$ws.Cells.Item(138, 5).Value = 1335
$ws.Cells.Item(138, 6).Value = 0
$ws.Cells.Item(138, 7).Value = 0
$ws.Cells.Item(138, 8).Value = 50

$ws.Cells.Item(139, 1).Value = "Benin"
$ws.Cells.Item(139, 2).Value = 2145
$ws.Cells.Item(139, 3).Value = 0
$ws.Cells.Item(139, 4).Value = 1738
$ws.Cells.Item(139, 5).Value = 367
$ws.Cells.Item(139, 6).Value = 0
$ws.Cells.Item(139, 7).Value = 0
$ws.Cells.Item(139, 8).Value = 40
